# Apply "ОК" marks to a handful of homework cells on the attendance/grades
# sheet, and move the active selection to E8 (which also re-homes the
# frozen-pane scroll position).
#
# Style note: several of the target cells did not previously exist in the
# sheet (no <c> element at all), so simply assigning .Value to them would
# leave them without the shared "centered/bordered" cell style (s="2") used
# by every other data cell in the table. To keep their formatting consistent
# with their row, we first copy the format from a neighboring cell in the
# same row (column C, which always carries style s="2") via
# Copy + PasteSpecial(xlPasteFormats), and only then set the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-MarkWithFormat($cellRef, $formatSourceRef) {
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats)
    $ws.Range($cellRef).Value = "ОК"
}

# --- Row 8 (Воробьева Полина) ---
$ws.Range("D8").Value = "ОК"
$ws.Range("G8").Value = "ОК"
Set-MarkWithFormat "I8" "C8"

# --- Row 9 (Горбенко Полина) ---
$ws.Range("E9").Value = "ОК"
$ws.Range("F9").Value = "ОК"
Set-MarkWithFormat "G9" "C9"

# --- Row 10 (Гришин Андрей) ---
Set-MarkWithFormat "I10" "C10"

# --- Row 13 (Емельяненко Семён) ---
$ws.Range("D13").Value = "ОК"
Set-MarkWithFormat "G13" "C13"
Set-MarkWithFormat "H13" "C13"

# --- Row 16 (Казаков Егор) ---
$ws.Range("D16").Value = "ОК"
$ws.Range("E16").Value = "ОК"
Set-MarkWithFormat "H16" "C16"

# --- Row 18 (Катахова Марина) ---
$ws.Range("D18").Value = "ОК"
$ws.Range("E18").Value = "ОК"
$ws.Range("F18").Value = "ОК"
Set-MarkWithFormat "H18" "C18"
Set-MarkWithFormat "I18" "C18"

# --- Row 30 (Фартушняк Василий) ---
$ws.Range("C30").Value = "ОК"

# Move the active selection to E8 (also updates the frozen pane's
# top-left-cell bookkeeping).
$ws.Range("E8").Select() | Out-Null
